$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.048.27"
$ws.Range("E2").Value = "  -3.25%  "

$ws.Range("D3").Value = "3.027.47"
$ws.Range("E3").Value = "  -3.85%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.14"
$ws.Range("E5").Value = "  -5.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.36"
$ws.Range("E6").Value = "  -8.03%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "3.023.37"
$ws.Range("E8").Value = "  -3.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -0.29%  "

$ws.Range("E10").Value = "  -2.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.05"
$ws.Range("E11").Value = "  -9.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.443"
$ws.Range("E12").Value = "  -4.03%  "

$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.35"
$ws.Range("E14").Value = "  -7.87%  "

$ws.Range("D15").Value = "3.514.37"
$ws.Range("E15").Value = "  -3.84%  "

$ws.Range("D16").Value = "61.973.51"
$ws.Range("E16").Value = "  -3.38%  "

$ws.Range("E17").Value = "  -2.40%  "

$ws.Range("D18").Value = "3.022.40"
$ws.Range("E18").Value = "  -4.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("E19").Value = "  -4.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "472.43"
$ws.Range("E20").Value = "  -7.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.01"
$ws.Range("E21").Value = "  -6.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("E22").Value = "  -4.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.88"
$ws.Range("E23").Value = "  -7.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.86"
$ws.Range("E24").Value = "  -1.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.74"
$ws.Range("E25").Value = "  -7.60%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.62"
$ws.Range("E27").Value = "  -6.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.96"
$ws.Range("E28").Value = "  -8.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.35"
$ws.Range("E30").Value = "  -4.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.82"
$ws.Range("E31").Value = "  -13.46%  "

$ws.Range("E32").Value = "  -4.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "56.46"
$ws.Range("E33").Value = "  +4.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.32"
$ws.Range("E34").Value = "  -10.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.18"
$ws.Range("E35").Value = "  -2.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.80"
$ws.Range("E36").Value = "  -4.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "464.97"
$ws.Range("E37").Value = "  -16.11%  "

$ws.Range("D38").Value = "3.058.17"
$ws.Range("E38").Value = "  -2.94%  "

$ws.Range("E39").Value = "  -9.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0775"
$ws.Range("E40").Value = "  -4.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("E41").Value = "  -7.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.91"
$ws.Range("E42").Value = "  -3.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.50"
$ws.Range("E43").Value = "  -7.79%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.243"
$ws.Range("E45").Value = "  -7.19%  "

$ws.Range("D46").Value = "0.0₃0529"
$ws.Range("E46").Value = "  +3.46%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.97"
$ws.Range("E47").Value = "  -8.73%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.28"
$ws.Range("E48").Value = "  -3.36%  "

$ws.Range("E49").Value = "  -1.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.80"
$ws.Range("E50").Value = "  -4.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.28"
$ws.Range("E51").Value = "  +2.86%  "
